$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.990.32"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "1.600.37"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'211.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D10").Value = "'18.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("D12").Value = "1.823.58"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "1.607.27"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "'4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "'0.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "25.995.00"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "'201.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.30%  "
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'6.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("E24").Value = "  +7.68%  "
$ws.Range("D25").Value = "'141.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E27").Value = "  -7.12%  "
$ws.Range("D28").Value = "'15.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "'6.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").Value = "'1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").Value = "1.125.43"
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("E37").Value = "  +10.00%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'0.786"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.735.64"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").Value = "'93.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").Value = "'53.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  -0.37%  "
